$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '30.699.62'
Set-TextValue $ws.Range('E2') '  +1.66%  '

Set-TextValue $ws.Range('D3') '1.897.74'
Set-TextValue $ws.Range('E3') '  +2.57%  '

Set-TextValue $ws.Range('E4') '  +0.15%  '

Set-TextValue $ws.Range('D5') '239.29'
Set-TextValue $ws.Range('E5') '  +1.28%  '

Set-TextValue $ws.Range('D6') '1.000'
Set-TextValue $ws.Range('E6') '  +0.07%  '

Set-TextValue $ws.Range('D7') '0.4829'

Set-TextValue $ws.Range('D8') '0.2854'
Set-TextValue $ws.Range('E8') '  +1.73%  '

Set-TextValue $ws.Range('D9') '0.06557'
Set-TextValue $ws.Range('E9') '  +1.27%  '

Set-TextValue $ws.Range('D10') '2.001.20'
Set-TextValue $ws.Range('E10') '  +8.06%  '

Set-TextValue $ws.Range('D11') '0.07464'
Set-TextValue $ws.Range('E11') '  +2.01%  '

Set-TextValue $ws.Range('D12') '16.75'
Set-TextValue $ws.Range('E12') '  +2.94%  '

Set-TextValue $ws.Range('D13') '5.111'
Set-TextValue $ws.Range('E13') '  +0.08%  '

Set-TextValue $ws.Range('D14') '88.07'
Set-TextValue $ws.Range('E14') '  +1.07%  '

Set-TextValue $ws.Range('D15') '0.6677'
Set-TextValue $ws.Range('E15') '  +3.38%  '

Set-TextValue $ws.Range('D16') '30.685.39'
Set-TextValue $ws.Range('E16') '  +1.86%  '

Set-TextValue $ws.Range('B17') 'Avalanche'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D17') '13.33'
Set-TextValue $ws.Range('E17') '  +0.77%  '

Set-TextValue $ws.Range('B18') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D18') '2.252.47'
Set-TextValue $ws.Range('E18') '  +7.40%  '

Set-TextValue $ws.Range('D19') '1.001'
Set-TextValue $ws.Range('E19') '  +0.11%  '

Set-TextValue $ws.Range('D20') '0.000007612'
Set-TextValue $ws.Range('E20') '  -0.16%  '

Set-TextValue $ws.Range('D21') '231.41'
Set-TextValue $ws.Range('E21') '  +2.58%  '

Set-TextValue $ws.Range('B22') 'Uniswap'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D22') '5.291'

Set-TextValue $ws.Range('B23') 'BinanceUSD'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D23') '1.001'
Set-TextValue $ws.Range('E23') '  +0.13%  '

Set-TextValue $ws.Range('D24') '6.245'
Set-TextValue $ws.Range('E24') '  +2.79%  '

Set-TextValue $ws.Range('D25') '169.93'
Set-TextValue $ws.Range('E25') '  +3.99%  '

Set-TextValue $ws.Range('D26') '9.343'
Set-TextValue $ws.Range('E26') '  +1.37%  '

Set-TextValue $ws.Range('D27') '18.78'
Set-TextValue $ws.Range('E27') '  +1.50%  '

Set-TextValue $ws.Range('E28') '  +2.75%  '

Set-TextValue $ws.Range('E29') '  -1.55%  '

Set-TextValue $ws.Range('E30') '  +11.00%  '

Set-TextValue $ws.Range('D31') '4.357'
Set-TextValue $ws.Range('E31') '  +2.80%  '

Set-TextValue $ws.Range('D32') '4.032'
Set-TextValue $ws.Range('E32') '  +1.98%  '

Set-TextValue $ws.Range('E33') '  +2.10%  '

Set-TextValue $ws.Range('E34') '  +6.90%  '

Set-TextValue $ws.Range('D35') '0.7598'
Set-TextValue $ws.Range('E35') '  +2.89%  '

Set-TextValue $ws.Range('D36') '2.708'
Set-TextValue $ws.Range('E36') '  +0.83%  '

Set-TextValue $ws.Range('E37') '  +4.27%  '

Set-TextValue $ws.Range('D38') '2.661'
Set-TextValue $ws.Range('E38') '  +2.31%  '

Set-TextValue $ws.Range('D39') '0.9216'
Set-TextValue $ws.Range('E39') '  +1.73%  '

Set-TextValue $ws.Range('D40') '2.081'
Set-TextValue $ws.Range('E40') '  +1.29%  '

Set-TextValue $ws.Range('D41') '107.10'
Set-TextValue $ws.Range('E41') '  +0.48%  '

Set-TextValue $ws.Range('D42') '0.4307'
Set-TextValue $ws.Range('E42') '  +1.36%  '

Set-TextValue $ws.Range('E43') '  +0.60%  '

Set-TextValue $ws.Range('D44') '5.727'
Set-TextValue $ws.Range('E44') '  -3.94%  '

Set-TextValue $ws.Range('D45') '7.444'
Set-TextValue $ws.Range('E45') '  +0.85%  '

Set-TextValue $ws.Range('D46') '64.65'
Set-TextValue $ws.Range('E46') '  +0.95%  '

Set-TextValue $ws.Range('E47') '  -3.16%  '

Set-TextValue $ws.Range('D48') '1.495'
Set-TextValue $ws.Range('E48') '  -4.05%  '

Set-TextValue $ws.Range('D49') '8.954'
Set-TextValue $ws.Range('E49') '  +2.76%  '

Set-TextValue $ws.Range('D50') '33.93'
Set-TextValue $ws.Range('E50') '  -0.65%  '

Set-TextValue $ws.Range('D51') '0.05681'
Set-TextValue $ws.Range('E51') '  +0.43%  '
